# Delete row 13 ("dat" / "dfyarchicloud" entry) from the "groups" worksheet.
# This shifts row 14 (1309 / ansible-roles) up to become the new row 13,
# shrinking the used range from A1:E14 to A1:E13.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("groups")
$ws.Rows.Item(13).Delete()

# Update the active-cell selections on each sheet to match the saved file.
$ws.Range("B13").Select()

$ws2 = $wb.Worksheets.Item("projects")
$ws2.Range("B18").Select()

$ws.Activate()
